# Rename the sheet from "moret_nucs_rating_no1_alpha1" to "S10" and keep
# the Print_Area / Print_Titles defined names (and their sheet-qualified
# references) pointing at the renamed sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "S10"

foreach ($n in $wb.Names) {
    if ($n.Name -eq "S10!Print_Area") {
        $n.RefersTo = "='S10'!`$A`$1:`$D`$470"
    }
    elseif ($n.Name -eq "S10!Print_Titles") {
        $n.RefersTo = "='S10'!`$1:`$1"
    }
}

# Move the selection / active cell to A23 (and drop any scrolled
# top-left-cell position so the view resets to the top of the sheet).
$ws.Range("A23").Select()
